# Weekly update: insert a new "Espinaca" price record for
# Vega Central Mapocho de Santiago at the top of the data block (row 410),
# pushing the existing rows 410-436 down to 411-437.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 410 (shifts rows 410:436 -> 411:437)
$ws.Rows("410:410").Insert()

# Populate the newly inserted row with the new week's record
$ws.Range("A410").Value = 9
$ws.Range("B410").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C410").Value = "Metropolitana"
$ws.Range("D410").Value = 44714
$ws.Range("E410").Value = 13
$ws.Range("F410").Value = 100112012
$ws.Range("G410").Value = "Espinaca"
$ws.Range("H410").Value = "Sin especificar"
$ws.Range("I410").Value = "Primera"
$ws.Range("J410").Value = 160
$ws.Range("K410").Value = 6000
$ws.Range("L410").Value = 7000
$ws.Range("M410").Value = 6500
$ws.Range("N410").Value = "`$/cuna 10 kilos"
$ws.Range("O410").Value = "Provincia de Chacabuco"
$ws.Range("P410").Value = 650
$ws.Range("Q410").Value = 10
$ws.Range("R410").Value = "Hortaliza"
